$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2.013391040896522
$ws.Range("C2").Value = 0.4545716620954181
$ws.Range("D2").Value = 0.0548242614504082
$ws.Range("F2").Value = 3.224396755779395
$ws.Range("G2").Value = 0.002543731212099035
$ws.Range("J2").Value = 0.1999677489458449
$ws.Range("M2").Value = 0.6547522193246067
$ws.Range("N2").Value = 2.13922535829785
$ws.Range("B3").Value = 1.894004493856016
$ws.Range("C3").Value = 0.4211297235562483
$ws.Range("D3").Value = 0.05472724908364168
$ws.Range("F3").Value = 3.180009302102007
$ws.Range("G3").Value = 0.002549823428066786
$ws.Range("J3").Value = 0.199772679028726
$ws.Range("M3").Value = 0.6291005810636818
$ws.Range("N3").Value = 2.151440602265737
$ws.Range("B4").Value = 1.821991698144018
$ws.Range("C4").Value = 0.4008820743332251
$ws.Range("D4").Value = 0.05467032729818411
$ws.Range("F4").Value = 3.154611315949396
$ws.Range("G4").Value = 0.002553758468885161
$ws.Range("J4").Value = 0.1997365924836103
$ws.Range("M4").Value = 0.6137668284747875
$ws.Range("N4").Value = 2.159624172293071
$ws.Range("B5").Value = 1.792967871958126
$ws.Range("C5").Value = 0.392702009086122
$ws.Range("D5").Value = 0.05464779058048208
$ws.Range("F5").Value = 3.144725662258395
$ws.Range("G5").Value = 0.002555411089922224
$ws.Range("J5").Value = 0.1997429023989454
$ws.Range("M5").Value = 0.607622376754037
$ws.Range("N5").Value = 2.163130167556652
$ws.Range("B6").Value = 1.788167857204655
$ws.Range("C6").Value = 0.3913479798583523
$ws.Range("D6").Value = 0.05464408804617493
$ws.Range("F6").Value = 3.143112123137996
$ws.Range("G6").Value = 0.002555688474813567
$ws.Range("J6").Value = 0.1997452186317759
$ws.Range("M6").Value = 0.6066083716073862
$ws.Range("N6").Value = 2.163722650142475
$ws.Range("B7").Value = 1.821598972499771
$ws.Range("C7").Value = 0.4007714686961492
$ws.Range("D7").Value = 0.0546700206970544
$ws.Range("F7").Value = 3.154476118170606
$ws.Range("G7").Value = 0.002553780557831358
$ws.Range("J7").Value = 0.1997365925266266
$ws.Range("M7").Value = 0.613683541149399
$ws.Range("N7").Value = 2.159670763426931
$ws.Range("B8").Value = 1.971956947317949
$ws.Range("C8").Value = 0.4429810241281871
$ws.Range("D8").Value = 0.05479025930429438
$ws.Range("F8").Value = 3.208705196551989
$ws.Range("G8").Value = 0.002545791568234805
$ws.Range("J8").Value = 0.1998831017278135
$ws.Range("M8").Value = 0.6458207580056268
$ws.Range("N8").Value = 2.143294912455517
$ws.Range("B9").Value = 2.277191495271154
$ws.Range("C9").Value = 0.5280656699717383
$ws.Range("D9").Value = 0.05504731160352505
$ws.Range("F9").Value = 3.329898514897252
$ws.Range("G9").Value = 0.002531659508293089
$ws.Range("J9").Value = 0.2008359734125307
$ws.Range("M9").Value = 0.7121752512898354
$ws.Range("N9").Value = 2.116636350198547
$ws.Range("B10").Value = 2.507997490032778
$ws.Range("C10").Value = 0.5920547887501471
$ws.Range("D10").Value = 0.0552495756243252
$ws.Range("F10").Value = 3.42817777014352
$ws.Range("G10").Value = 0.002522200658512144
$ws.Range("J10").Value = 0.2019443304191171
$ws.Range("M10").Value = 0.7630049846874556
$ws.Range("N10").Value = 2.100420355881795
$ws.Range("B11").Value = 2.614469513131496
$ws.Range("C11").Value = 0.6215012849884829
$ws.Range("D11").Value = 0.05534460768776217
$ws.Range("F11").Value = 3.47493519891313
$ws.Range("G11").Value = 0.002518095777017925
$ws.Range("J11").Value = 0.2025377975353209
$ws.Range("M11").Value = 0.7865909407960601
$ws.Range("N11").Value = 2.093784926428228
$ws.Range("B12").Value = 2.655003745865372
$ws.Range("C12").Value = 0.6327015756172614
$ws.Range("D12").Value = 0.05538103700890162
$ws.Range("F12").Value = 3.492938951279285
$ws.Range("G12").Value = 0.002516569649551062
$ws.Range("J12").Value = 0.2027754092559348
$ws.Range("M12").Value = 0.7955897148610944
$ws.Range("N12").Value = 2.091379698657875
$ws.Range("B13").Value = 2.646264334059367
$ws.Range("C13").Value = 0.6302871730186439
$ws.Range("D13").Value = 0.05537317147492415
$ws.Range("F13").Value = 3.489048231393895
$ws.Range("G13").Value = 0.002516897072370334
$ws.Range("J13").Value = 0.2027236618620734
$ws.Range("M13").Value = 0.7936486645299681
$ws.Range("N13").Value = 2.091892914000354
$ws.Range("B14").Value = 2.617799949750975
$ws.Range("C14").Value = 0.6224217409150015
$ws.Range("D14").Value = 0.05534759582725535
$ws.Range("F14").Value = 3.476410394240219
$ws.Range("G14").Value = 0.002517969655270969
$ws.Range("J14").Value = 0.2025570876445641
$ws.Range("M14").Value = 0.7873299229817974
$ws.Range("N14").Value = 2.093584888120276
$ws.Range("B15").Value = 2.600392863628883
$ws.Range("C15").Value = 0.6176104186321822
$ws.Range("D15").Value = 0.05533198791774963
$ws.Range("F15").Value = 3.468708216653567
$ws.Range("G15").Value = 0.002518630323934583
$ws.Range("J15").Value = 0.2024567345623893
$ws.Range("M15").Value = 0.7834682938936623
$ws.Range("N15").Value = 2.094635292209745
$ws.Range("B16").Value = 2.501069287353118
$ws.Range("C16").Value = 0.5901372683450745
$ws.Range("D16").Value = 0.05524342654921721
$ws.Range("F16").Value = 3.425163561314974
$ws.Range("G16").Value = 0.002522472891843932
$ws.Range("J16").Value = 0.201907345818654
$ws.Range("M16").Value = 0.7614729656436765
$ws.Range("N16").Value = 2.10086898755614
$ws.Range("B17").Value = 2.440518071673409
$ws.Range("C17").Value = 0.5733705378793843
$ws.Range("D17").Value = 0.05518987657272167
$ws.Range("F17").Value = 3.398977409752348
$ws.Range("G17").Value = 0.002524880773714032
$ws.Range("J17").Value = 0.2015932070125501
$ws.Range("M17").Value = 0.7480986653817752
$ws.Range("N17").Value = 2.104883636630902
$ws.Range("B18").Value = 2.405829463066027
$ws.Range("C18").Value = 0.5637584660235007
$ws.Range("D18").Value = 0.05515935988132092
$ws.Range("F18").Value = 3.384108573999271
$ws.Range("G18").Value = 0.002526284370176102
$ws.Range("J18").Value = 0.2014209226129893
$ws.Range("M18").Value = 0.7404496843532797
$ws.Range("N18").Value = 2.107262479430389
$ws.Range("B19").Value = 2.394108248382565
$ws.Range("C19").Value = 0.5605093992342631
$ws.Range("D19").Value = 0.05514907595836016
$ws.Range("F19").Value = 3.379107255943012
$ws.Range("G19").Value = 0.002526762811454221
$ws.Range("J19").Value = 0.2013640315244842
$ws.Range("M19").Value = 0.73786733265203
$ws.Range("N19").Value = 2.108079863490261
$ws.Range("B20").Value = 2.446949458888753
$ws.Range("C20").Value = 0.5751520960790231
$ws.Range("D20").Value = 0.0551955476242294
$ws.Range("F20").Value = 3.40174500037034
$ws.Range("G20").Value = 0.002524622521862196
$ws.Range("J20").Value = 0.201625777966342
$ws.Range("M20").Value = 0.74951786910934
$ws.Range("N20").Value = 2.104449048104641
$ws.Range("B21").Value = 2.626154752813306
$ws.Range("C21").Value = 0.6247306563194002
$ws.Range("D21").Value = 0.05535509593069921
$ws.Range("F21").Value = 3.480114329025696
$ws.Range("G21").Value = 0.002517653844948052
$ws.Range("J21").Value = 0.2026056646690293
$ws.Range("M21").Value = 0.7891840590474146
$ws.Range("N21").Value = 2.093084990758186
$ws.Range("B22").Value = 2.74453390664064
$ws.Range("C22").Value = 0.6574223596741717
$ws.Range("D22").Value = 0.05546195479217242
$ws.Range("F22").Value = 3.533070031476399
$ws.Range("G22").Value = 0.002513264300321357
$ws.Range("J22").Value = 0.2033211640025243
$ws.Range("M22").Value = 0.815500714958219
$ws.Range("N22").Value = 2.086284723976306
$ws.Range("B23").Value = 2.681236540002658
$ws.Range("C23").Value = 0.6399473782400378
$ws.Range("D23").Value = 0.05540468283349376
$ws.Range("F23").Value = 3.504646657351515
$ws.Range("G23").Value = 0.00251559205138177
$ws.Range("J23").Value = 0.2029324039523388
$ws.Range("M23").Value = 0.8014188803649915
$ws.Range("N23").Value = 2.089856512316473
$ws.Range("B24").Value = 2.44404144556313
$ws.Range("C24").Value = 0.5743465685197293
$ws.Range("D24").Value = 0.05519298290176011
$ws.Range("F24").Value = 3.40049319374134
$ws.Range("G24").Value = 0.002524739217460048
$ws.Range("J24").Value = 0.2016110267320315
$ws.Range("M24").Value = 0.7488761223067399
$ws.Range("N24").Value = 2.104645305217772
$ws.Range("B25").Value = 2.193483348231837
$ws.Range("C25").Value = 0.504793554893979
$ws.Range("D25").Value = 0.05497545129032089
$ws.Range("F25").Value = 3.295503705877252
$ws.Range("G25").Value = 0.002535319524473712
$ws.Range("J25").Value = 0.2005066856343021
$ws.Range("M25").Value = 0.6938629398850296
$ws.Range("N25").Value = 2.123259667674205
